# Add support for a second file (File name 2 / File use 2) to the
# "valid_columns" list, inserted right after the existing "File use" row.
#
# Original layout (rows 1-6):
#   1 Object Unique ID
#   2 Level
#   3 File name
#   4 File use
#   5 Type of Resource
#   6 Language
#
# New layout (rows 1-8):
#   1 Object Unique ID
#   2 Level
#   3 File name
#   4 File use
#   5 File name 2   <- new
#   6 File use 2    <- new
#   7 Type of Resource
#   8 Language

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 5 ("Type of Resource"),
# pushing everything else down. The inserted rows inherit the formatting
# of the row above (row 4), matching the style used by the rows around
# them in the target file.
$ws.Rows("5:6").Insert()

$ws.Range("A5").Value = "File name 2"
$ws.Range("A6").Value = "File use 2"

# Match the resulting selection left behind by this edit in Excel.
$ws.Rows("5:5").Select()
